# v2.0.0 - multiple items table[Hs
# Add a new "Shields" entry row to the bottom of the "models" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row right below the existing table (row 15 -> row 16),
# reusing the existing "Shields" category label already used on the
# "items" sheet (column A).
$ws.Range("A16").Value = "Shields"

# Move the selection to the next empty row, as left by the editor.
$ws.Range("A17").Select()
